$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new transaction was logged; insert a row above the current top data
# row (row 35) which pushes every existing row (35-92) down by one,
# stretching the sheet's used range from A1:Y92 to A1:Y93.
$ws.Rows.Item(35).Insert()

# Populate the freshly inserted row with the new entry.
$ws.Range("R35").Value = "balance your axis"
$ws.Range("S35").Value = "2024-09-08 09:53:37"
